$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of Spanish connector words (de/el/las -> De/El/Las)
$ws.Range("B6").Value = "Amatenango De La Frontera"
$ws.Range("B7").Value = "Comitán De Domínguez"
$ws.Range("B16").Value = "San Cristóbal De Las Casas"
$ws.Range("A20").Value = "Ciudad De México"
$ws.Range("A25").Value = "Estado De México"
$ws.Range("B27").Value = "Apaseo El Alto"
$ws.Range("B32").Value = "Chilapa De Álvarez"
$ws.Range("B38").Value = "Lagos De Moreno"
$ws.Range("B39").Value = "San Miguel El Alto"
$ws.Range("B41").Value = "Huajuapan De León"
$ws.Range("B42").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B43").Value = "Oaxaca De Juárez"
$ws.Range("B49").Value = "Chalchicomula De Sesma"
$ws.Range("A70").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B73").Value = "Martínez De La Torre"

# Fix the TOTAL row label casing
$ws.Range("A78").Value = "Total"

# Remove the trailing footnote/source rows (80-84)
$ws.Range("A80:D84").EntireRow.Delete()
